$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price value "looks numeric" (e.g. "246.42") need to be
# forced to stay stored as text (matching the original inlineStr cells),
# otherwise Excel auto-converts them to a numeric cell (losing trailing
# zeros / exact formatting and changing the cell type). Force text via a
# temporary "@" number format, then restore a clean "Normal" style so we
# don't leave a custom number format attached to the cell.
$numericLookingPriceCells = @(
    "D5","D6","D8","D9","D10","D11","D12","D14","D15","D17","D19",
    "D21","D22","D23","D25","D26","D27","D29","D32","D33","D36","D37",
    "D38","D40","D41","D42","D43","D46","D49","D50"
)
foreach ($addr in $numericLookingPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "35.307.93"
$ws.Range("E2").Value = "  -0.53%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.900.30"
$ws.Range("E3").Value = "  -0.15%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.10%  "

# Row 5 - was XRP, now BNB
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "246.42"
$ws.Range("E5").Value = "  +0.33%  "

# Row 6 - was BNB, now XRP
$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").Value = "0.692"
$ws.Range("E6").Value = "  +9.53%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.16%  "

# Row 8 - Solana
$ws.Range("D8").Value = "40.49"
$ws.Range("E8").Value = "  -3.32%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.351"
$ws.Range("E9").Value = "  +3.83%  "

# Row 10 - OKB
$ws.Range("D10").Value = "52.04"
$ws.Range("E10").Value = "  +6.86%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "0.0723"
$ws.Range("E11").Value = "  +2.78%  "

# Row 12 - TRON
$ws.Range("D12").Value = "0.0988"
$ws.Range("E12").Value = "  -0.86%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("E13").Value = "  -0.37%  "

# Row 14 - Chainlink
$ws.Range("D14").Value = "12.53"
$ws.Range("E14").Value = "  +1.87%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "0.710"
$ws.Range("E15").Value = "  +2.99%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "1.897.60"
$ws.Range("E16").Value = "  -0.51%  "

# Row 17 - Polkadot
$ws.Range("D17").Value = "4.85"
$ws.Range("E17").Value = "  +0.06%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "35.292.94"
$ws.Range("E18").Value = "  -0.55%  "

# Row 19 - Litecoin
$ws.Range("D19").Value = "72.54"
$ws.Range("E19").Value = "  +0.80%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.0₃0820"
$ws.Range("E20").Value = "  -0.10%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "241.16"
$ws.Range("E21").Value = "  -0.77%  "

# Row 22 - Avalanche
$ws.Range("D22").Value = "12.83"
$ws.Range("E22").Value = "  +2.75%  "

# Row 23 - Uniswap
$ws.Range("D23").Value = "5.06"
$ws.Range("E23").Value = "  +4.32%  "

# Row 24 - Dai
$ws.Range("E24").Value = "  -0.09%  "

# Row 25 - Toncoin
$ws.Range("D25").Value = "2.32"
$ws.Range("E25").Value = "  +1.21%  "

# Row 26 - PancakeSwap
$ws.Range("D26").Value = "2.30"
$ws.Range("E26").Value = "  +5.52%  "

# Row 27 - Monero
$ws.Range("D27").Value = "168.14"
$ws.Range("E27").Value = "  -2.20%  "

# Row 28 - Cosmos
$ws.Range("E28").Value = "  +0.36%  "

# Row 29 - EthereumClassic
$ws.Range("D29").Value = "18.97"
$ws.Range("E29").Value = "  +5.87%  "

# Row 30 - Stellar
$ws.Range("E30").Value = "  +4.65%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "4.23"
$ws.Range("E32").Value = "  +3.32%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "0.0572"
$ws.Range("E33").Value = "  +0.72%  "

# Row 34 - BinanceUSD
$ws.Range("E34").Value = "  -0.16%  "

# Row 35 - WEMIXToken
$ws.Range("E35").Value = "  +6.89%  "

# Row 36 - InternetComputer(DFINITY)
$ws.Range("D36").Value = "4.16"
$ws.Range("E36").Value = "  -0.28%  "

# Row 37 - ImmutableX
$ws.Range("D37").Value = "0.913"
$ws.Range("E37").Value = "  -5.21%  "

# Row 38 - TrustWalletToken
$ws.Range("D38").Value = "1.47"
$ws.Range("E38").Value = "  +7.06%  "

# Row 39 - LidoDAOToken
$ws.Range("E39").Value = "  -0.18%  "

# Row 40 - Aave
$ws.Range("D40").Value = "95.98"
$ws.Range("E40").Value = "  +5.83%  "

# Row 41 - ARBITRUM
$ws.Range("D41").Value = "1.10"
$ws.Range("E41").Value = "  -0.64%  "

# Row 42 - Kaspa
$ws.Range("D42").Value = "0.0650"
$ws.Range("E42").Value = "  +7.22%  "

# Row 43 - InjectiveProtocol
$ws.Range("D43").Value = "16.40"
$ws.Range("E43").Value = "  +3.98%  "

# Row 44 - VeChain
$ws.Range("E44").Value = "  +1.51%  "

# Row 45 - Maker
$ws.Range("D45").Value = "1.354.43"
$ws.Range("E45").Value = "  +0.23%  "

# Row 46 - RenderToken
$ws.Range("D46").Value = "2.40"
$ws.Range("E46").Value = "  +2.29%  "

# Row 47 - HuobiToken
$ws.Range("E47").Value = "  -0.08%  "

# Row 48 - MXToken
$ws.Range("E48").Value = "  +1.02%  "

# Row 49 - was MultiversX, now Gas
$ws.Range("B49").Value = "Gas"
$ws.Range("C49").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D49").Value = "12.45"
$ws.Range("E49").Value = "  -2.82%  "

# Row 50 - was Gas, now MultiversX
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").Value = "45.58"
$ws.Range("E50").Value = "  -6.81%  "

# Row 51 - FraxShare
$ws.Range("E51").Value = "  -2.27%  "

# Restore a plain "Normal" style on the cells we temporarily reformatted
# above so no stray number format remains attached to them.
foreach ($addr in $numericLookingPriceCells) {
    $ws.Range($addr).Style = "Normal"
}
